$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (A11:F13 "5/Set different Screen Size") was the last data row;
# add a new row 14 describing the "Save options in database" task.
# Row 12 already has the exact style pattern we need
# (A-D = centered/bordered, E-F = centered/bordered/wrap), so copy its
# formats down into row 14 before writing the new values.
$ws.Range("A12:F12").Copy()
$ws.Range("A14:F14").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A14").Value = 6
$ws.Range("B14").Value = "Save options in database"
$ws.Range("C14").Value = "When the application restarted, save the options and keep them"
$ws.Range("D14").Value = "Yes"
$ws.Range("E14").Value = "Create new table for options (Text Size, Text Color, and Background Color), `nthen call it on onCreate of Main Activity LifeCycle, and adapt the options."
$ws.Range("F14").Value = "MainActivity`nDBHelper`nAddNewWordDialog`nSettingDialog"

# Matches the row height recorded for the new row in the authored workbook.
$ws.Rows.Item(14).RowHeight = 57.6

# Move the selection down to reflect where the user ended up after typing
# the new row (one blank row below the newly added data).
[void]$ws.Range("D15").Select()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 4
